# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2404"
#   "<header>_new" -> "<header>_FV2410"
# then turn the data range into a real Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row (row 1) cells ---------------------------------
$lastCol = $ws.UsedRange.Columns.Count
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = ($val -replace "_old$", "_FV2404")
        } elseif ($val -like "*_new") {
            $cell.Value = ($val -replace "_new$", "_FV2410")
        }
    }
}

# --- 2) Turn the used range into an Excel Table ------------------------------
$lastRow = $ws.UsedRange.Rows.Count
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row (split below row 1) ----------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header columns renamed, Table1 created, header row frozen."
